# Generate Report for Handback
# Updates the status of the 1a2c766d-8b1c-4062-82fa-36ca7862d7ee.md file from
# "Ready for handoff" to "Handback transform failed" and records the handback
# transform error details for both the zh-cn and de-de locales.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update status text on the Overview sheet (row 3 = 1a2c766d...md)
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"

# Update status text on the per-locale sheets (row 3 = 1a2c766d...md)
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# Record the "Error Detail" column (P) values describing the handback mismatch
$zhcn.Range("P3").Value = "Handback file name: ptid2kdy.vqx is different with handoff file name: 1a2c766d-8b1c-4062-82fa-36ca7862d7ee.07ac5614b163c6659063d05d00da37cd89deb28c.zh-cn."
$dede.Range("P3").Value = "Handback file name: ptid2kdy.vqx is different with handoff file name: 1a2c766d-8b1c-4062-82fa-36ca7862d7ee.07ac5614b163c6659063d05d00da37cd89deb28c.de-de."

# Widen the "Error Detail" column (P, the 16th column) so the new messages are readable.
# ColumnWidth values are persisted with an extra 5/6 padding added on save, so back that
# off here to land on an effective stored width of exactly 40.
$targetColumnWidth = 40 - (5 / 6)
$zhcn.Columns.Item(16).ColumnWidth = $targetColumnWidth
$dede.Columns.Item(16).ColumnWidth = $targetColumnWidth
